$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# C10 currently holds 18 (the "From" hour for rule R40 / "Good Night").
# Per the target revision it should hold 1 instead.
$ws.Range("C10").Value = 1
